$wb = $excel.ActiveWorkbook

# --- DataTable1: add "notes:" / "this is not collected by ExcellBDD" row ---
$ws1 = $wb.Worksheets.Item("DataTable1")
$ws1.Activate()
$ws1.Range("A10").Value = "notes:"
$ws1.Range("B10").Value = "this is not collected by ExcellBDD"
$ws1.Range("A10").Select()

# --- DataTable2: add combined notes cell ---
$ws2 = $wb.Worksheets.Item("DataTable2")
$ws2.Activate()
$ws2.Range("B9").Value = "notes: this is not collected by ExcellBDD"
$ws2.Range("B9").Select()

# --- DataTable3: add combined notes cell ---
$ws3 = $wb.Worksheets.Item("DataTable3")
$ws3.Activate()
$ws3.Range("C11").Value = "notes: this is not collected by ExcellBDD"
$ws3.Range("C11").Select()

# --- DataTableV0.5: add combined notes cell ---
$ws4 = $wb.Worksheets.Item("DataTableV0.5")
$ws4.Activate()
$ws4.Range("A10").Value = "notes: this is not collected by ExcellBDD"
$ws4.Range("A10").Select()

# --- Back to the main sheet, move selection ---
$wsBDD = $wb.Worksheets.Item("DataTableBDD")
$wsBDD.Activate()
$wsBDD.Range("D19").Select()

# --- Re-align the workbook window position ---
$wb.Windows.Item(1).Left = -120
